# Auto-generated edit script applying numeric corrections to H..N columns
# across multiple sheets, per the scheduled profit-data refresh.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2163.3635
$ws.Range("I28").Value = 982.8333
$ws.Range("J28").Value = 3580
$ws.Range("K28").Value = 982.8333
$ws.Range("L28").Value = 3580
$ws.Range("M28").Value = -497.8333
$ws.Range("N28").Value = -4550
$ws.Range("H32").Value = 5000
$ws.Range("J32").Value = 5000
$ws.Range("L32").Value = 5000
$ws.Range("N32").Value = -5652
$ws.Range("H55").Value = 856.4375
$ws.Range("I55").Value = 752.5
$ws.Range("J55").Value = 918.8
$ws.Range("K55").Value = 752.5
$ws.Range("L55").Value = 918.8
$ws.Range("M55").Value = -538.5
$ws.Range("N55").Value = -1346.8
$ws.Range("H62").Value = 13062
$ws.Range("I62").Value = 15129.375
$ws.Range("K62").Value = 15129.375
$ws.Range("M62").Value = -14505.375
$ws.Range("H65").Value = 13062
$ws.Range("I65").Value = 15129.375
$ws.Range("K65").Value = 75646.875
$ws.Range("M65").Value = -72526.875
$ws.Range("H94").Value = 1833.6
$ws.Range("I94").Value = 1833.6
$ws.Range("K94").Value = 1833.6
$ws.Range("M94").Value = -1382.6
$ws.Range("H111").Value = 2864.8333
$ws.Range("I111").Value = 2864.8333
$ws.Range("K111").Value = 8594.499899999999
$ws.Range("M111").Value = -5527.499899999999
$ws.Range("H118").Value = 1070.9
$ws.Range("J118").Value = 1195
$ws.Range("L118").Value = 3585
$ws.Range("N118").Value = -6899
$ws.Range("H138").Value = 2748.2856
$ws.Range("I138").Value = 2460.8157
$ws.Range("K138").Value = 7382.4471
$ws.Range("M138").Value = -2242.4471

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2617.318
$ws.Range("I74").Value = 2337.7878
$ws.Range("K74").Value = 2337.7878
$ws.Range("M74").Value = -1463.7878
$ws.Range("H77").Value = 2617.318
$ws.Range("I77").Value = 2337.7878
$ws.Range("K77").Value = 11688.939
$ws.Range("M77").Value = -7320.939
$ws.Range("H122").Value = 2124.375
$ws.Range("I122").Value = 1749.2858
$ws.Range("J122").Value = 4750
$ws.Range("K122").Value = 5247.857400000001
$ws.Range("L122").Value = 14250
$ws.Range("M122").Value = -2797.857400000001
$ws.Range("N122").Value = -19150
$ws.Range("H132").Value = 22906.98
$ws.Range("I132").Value = 26609.904
$ws.Range("J132").Value = 3466.625
$ws.Range("K132").Value = 79829.712
$ws.Range("L132").Value = 10399.875
$ws.Range("M132").Value = -77299.712
$ws.Range("N132").Value = -15459.875

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 6450.154
$ws.Range("I22").Value = 4585.7
$ws.Range("K22").Value = 4585.7
$ws.Range("M22").Value = -4412.7
$ws.Range("H26").Value = 42104
$ws.Range("I26").Value = 42632.332
$ws.Range("J26").Value = 40519
$ws.Range("K26").Value = 42632.332
$ws.Range("L26").Value = 40519
$ws.Range("M26").Value = -42340.332
$ws.Range("N26").Value = -41103
$ws.Range("H81").Value = 75999.5
$ws.Range("J81").Value = 75999.5
$ws.Range("L81").Value = 75999.5
$ws.Range("N81").Value = -78121.5
$ws.Range("H84").Value = 75999.5
$ws.Range("J84").Value = 75999.5
$ws.Range("L84").Value = 227998.5
$ws.Range("N84").Value = -238606.5
$ws.Range("H86").Value = 5015.1055
$ws.Range("I86").Value = 4937.1665
$ws.Range("J86").Value = 5051.077
$ws.Range("K86").Value = 4937.1665
$ws.Range("L86").Value = 5051.077
$ws.Range("M86").Value = -3814.1665
$ws.Range("N86").Value = -7297.077
$ws.Range("H89").Value = 5015.1055
$ws.Range("I89").Value = 4937.1665
$ws.Range("J89").Value = 5051.077
$ws.Range("K89").Value = 24685.8325
$ws.Range("L89").Value = 25255.385
$ws.Range("M89").Value = -19069.8325
$ws.Range("N89").Value = -36487.385
$ws.Range("H99").Value = 6256.909
$ws.Range("I99").Value = 5595.1113
$ws.Range("K99").Value = 5595.1113
$ws.Range("M99").Value = -4097.1113

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1067.4445
$ws.Range("I22").Value = 393.5
$ws.Range("J22").Value = 1909.875
$ws.Range("K22").Value = 393.5
$ws.Range("L22").Value = 1909.875
$ws.Range("M22").Value = -43.5
$ws.Range("N22").Value = -2609.875
$ws.Range("H58").Value = 49950.094
$ws.Range("I58").Value = 68706.664
$ws.Range("K58").Value = 68706.664
$ws.Range("M58").Value = -68503.664
$ws.Range("H132").Value = 2528.8333
$ws.Range("I132").Value = 2560.0908
$ws.Range("J132").Value = 2185
$ws.Range("K132").Value = 7680.2724
$ws.Range("L132").Value = 6555
$ws.Range("M132").Value = -5150.2724
$ws.Range("N132").Value = -11615
$ws.Range("H134").Value = 57010.473
$ws.Range("I134").Value = 58621.39
$ws.Range("J134").Value = 28014
$ws.Range("K134").Value = 175864.17
$ws.Range("L134").Value = 84042
$ws.Range("M134").Value = -173329.17
$ws.Range("N134").Value = -89112
$ws.Range("H136").Value = 49950.094
$ws.Range("I136").Value = 68706.664
$ws.Range("K136").Value = 206119.992
$ws.Range("M136").Value = -203569.992

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 160.5
$ws.Range("I12").Value = 84
$ws.Range("J12").Value = 206.4
$ws.Range("K12").Value = 252
$ws.Range("L12").Value = 619.2
$ws.Range("M12").Value = -79
$ws.Range("N12").Value = -965.2
$ws.Range("H55").Value = 7207.077
$ws.Range("J55").Value = 7732.6665
$ws.Range("L55").Value = 23197.9995
$ws.Range("N55").Value = -23551.9995
$ws.Range("H98").Value = 1406.5385
$ws.Range("I98").Value = 787.6667
$ws.Range("J98").Value = 1937
$ws.Range("K98").Value = 2363.0001
$ws.Range("L98").Value = 5811
$ws.Range("M98").Value = -865.0001000000002
$ws.Range("N98").Value = -8807
$ws.Range("H131").Value = 10896.091
$ws.Range("I131").Value = 1069.9
$ws.Range("J131").Value = 19084.584
$ws.Range("K131").Value = 3209.7
$ws.Range("L131").Value = 57253.75199999999
$ws.Range("M131").Value = 1830.3
$ws.Range("N131").Value = -67333.75199999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 20716.75
$ws.Range("I41").Value = 21438.5
$ws.Range("K41").Value = 21438.5
$ws.Range("M41").Value = -21083.5
$ws.Range("H122").Value = 3540.7646
$ws.Range("I122").Value = 2699.111
$ws.Range("K122").Value = 8097.333
$ws.Range("M122").Value = -5647.333
$ws.Range("H126").Value = 4926.75
$ws.Range("I126").Value = 3915.818
$ws.Range("J126").Value = 6162.3335
$ws.Range("K126").Value = 11747.454
$ws.Range("L126").Value = 18487.0005
$ws.Range("M126").Value = -9277.454000000002
$ws.Range("N126").Value = -23427.0005
$ws.Range("H132").Value = 69122.87
$ws.Range("I132").Value = 79180.38
$ws.Range("K132").Value = 237541.14
$ws.Range("M132").Value = -235011.14

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2852.9443
$ws.Range("I7").Value = 2010.5333
$ws.Range("K7").Value = 2010.5333
$ws.Range("M7").Value = -1898.5333
$ws.Range("H46").Value = 2239.8147
$ws.Range("J46").Value = 2268.2693
$ws.Range("L46").Value = 2268.2693
$ws.Range("N46").Value = -2644.2693
$ws.Range("H93").Value = 1485.5667
$ws.Range("I93").Value = 1489.875
$ws.Range("J93").Value = 1468.3334
$ws.Range("K93").Value = 1489.875
$ws.Range("L93").Value = 1468.3334
$ws.Range("M93").Value = -241.875
$ws.Range("N93").Value = -3964.3334
$ws.Range("H126").Value = 2852.9443
$ws.Range("I126").Value = 2010.5333
$ws.Range("K126").Value = 6031.5999
$ws.Range("M126").Value = -3561.5999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 18102.334
$ws.Range("J74").Value = 18763
$ws.Range("L74").Value = 18763
$ws.Range("N74").Value = -20635
$ws.Range("H77").Value = 18102.334
$ws.Range("J77").Value = 18763
$ws.Range("L77").Value = 56289
$ws.Range("N77").Value = -65649
$ws.Range("H107").Value = 2157.0588
$ws.Range("I107").Value = 1219.5555
$ws.Range("J107").Value = 3211.75
$ws.Range("K107").Value = 3658.6665
$ws.Range("L107").Value = 9635.25
$ws.Range("M107").Value = -1738.6665
$ws.Range("N107").Value = -13475.25
$ws.Range("H122").Value = 1995.6666
$ws.Range("I122").Value = 1995.6666
$ws.Range("K122").Value = 5986.9998
$ws.Range("M122").Value = -3536.9998
